# Edit script: Cashflow sensitivity data + report cleanup
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cashflow")
$main = $wb.Worksheets.Item("Main")

# --- Row 2 gains the trailing (spare) styled BO:BV cells that row 3/4 already had ---
$ws.Range("BO3:BV3").Copy($ws.Range("BO2:BV2"))

# --- Row 3: change AocType/Novelty from CL/C to AU/I ---
$ws.Range("D3").Value = "AU"
$ws.Range("E3").Value = "I"

# --- Row 4: change AmountType from PR to CL, AocType/Novelty to EV/I ---
$ws.Range("B4").Value = "CL"
$ws.Range("D4").Value = "EV"
$ws.Range("E4").Value = "I"

# --- Use row 4 (fully styled incl. spare BO:BV cells) as a style template
#     for the five new rows we are about to create ---
$ws.Range("A4:BV4").Copy($ws.Range("A5:BV5"))
$ws.Range("A4:BV4").Copy($ws.Range("A6:BV6"))
$ws.Range("A4:BV4").Copy($ws.Range("A7:BV7"))
$ws.Range("A4:BV4").Copy($ws.Range("A8:BV8"))
$ws.Range("A4:BV4").Copy($ws.Range("A9:BV9"))

# --- Row 5: GIC1 / CL / BE / CL / C (same combination former row 3 had) ---
$ws.Range("A5").Value = "GIC1"
$ws.Range("B5").Value = "CL"
$ws.Range("C5").Value = "BE"
$ws.Range("D5").Value = "CL"
$ws.Range("E5").Value = "C"

# --- Row 6: GIC1 / PR / BE / BOP / N ---
$ws.Range("A6").Value = "GIC1"
$ws.Range("B6").Value = "PR"
$ws.Range("C6").Value = "BE"
$ws.Range("D6").Value = "BOP"
$ws.Range("E6").Value = "N"

# --- Row 7: GIC1 / PR / BE / AU / I ---
$ws.Range("A7").Value = "GIC1"
$ws.Range("B7").Value = "PR"
$ws.Range("C7").Value = "BE"
$ws.Range("D7").Value = "AU"
$ws.Range("E7").Value = "I"

# --- Row 8: GIC1 / PR / BE / EV / I ---
$ws.Range("A8").Value = "GIC1"
$ws.Range("B8").Value = "PR"
$ws.Range("C8").Value = "BE"
$ws.Range("D8").Value = "EV"
$ws.Range("E8").Value = "I"

# --- Row 9: GIC1 / PR / BE / CL / C ---
$ws.Range("A9").Value = "GIC1"
$ws.Range("B9").Value = "PR"
$ws.Range("C9").Value = "BE"
$ws.Range("D9").Value = "CL"
$ws.Range("E9").Value = "C"

# --- Refresh the Values0..Values60 (F:BN) numbers on every data row ---
$row2vals = @(-80,-80.079999999999984,-80.240159999999989,-80.480880479999982,-80.802804001919981,-81.206818021929578,-81.694058930061161,-82.265917342571584,-82.924044681312154,-83.67036108344395,-84.507064694278384,-85.521149470609728,-86.718445563198259,-88.10594069220943,-89.691847624669194,-91.485684577162573,-93.498369637860151,-95.742330509168795,-98.423115763425528,-101.57265546785514,-105.22927106469793,-109.43844190728585,-114.25373335120644,-119.73791255206434,-125.9642840047717,-133.01828390903893,-141.13239922749028,-150.44713757650464,-161.12888434443647,-173.37467955461366,-187.41802859853738,-203.53597905801161,-222.05775315229064,-243.37529745491057,-266.2525754156722,-290.48155977849837,-315.75345547922774,-341.96099228400362,-368.97591067443989,-396.28012806434845,-423.62345690078848,-450.735358142439,-476.8780089147005,-501.67566537826497,-523.74939465490866,-542.08062346783038,-555.0905584310583,-561.75164513223103,-563.43690006762768,-560.05627866722193,-546.61492797920857,-510.53834273258076,-463.56881520118333,-407.01341974663899,-343.11231284641667,-275.51918721567256,-207.46594797340143,-143.77390194556719,-89.571140912088353,-46.576993274285947,-13.973097982285781)
$row3vals = @(-73.600000000000009,-73.804183056030055,-74.005988873511555,-74.269718768748874,-74.602103085421973,-75.006373842903812,-75.484843073551531,-76.039564981379542,-76.672583056089877,-77.386045313005297,-78.182267965644314,-79.162546254399487,-80.310070532763987,-81.632147903865302,-83.137014389372965,-84.833991199275289,-86.733618778901089,-88.847787214432557,-91.398523020198041,-94.38378902608703,-97.84048680705294,-101.81202532505596,-106.34926686943705,-111.51162183060511,-117.36833499037918,-124.00000655936267,-131.63924730391071,-140.40441871180292,-150.45251100323148,-161.96935505694907,-175.17511790862838,-190.33099133776699,-207.74734818051408,-227.79370894463409,-249.16177239767649,-271.76161553762893,-295.29571363611007,-319.68428445689688,-344.80562754757858,-370.13651614737176,-395.47166973057449,-420.55648958365617,-444.64998164767547,-467.44073541042076,-487.51460902938481,-503.94367963936435,-515.11768679476836,-519.98471567021829,-519.95330893972323,-517.48879736462914,-507.25239871488412,-476.45864451893556,-433.73385895599188,-381.78314678385362,-322.67382619765033,-259.78356161859261,-196.18556690700947,-136.3798155198939,-85.241419067552911,-44.514839598436339,-13.458093512555834)
$row4vals = @(-80,-80.221938104380499,-80.441292253816911,-80.727955183422679,-81.089242484154312,-81.528667220547618,-82.048742471251657,-82.651701066716882,-83.33976419140204,-84.115266644570966,-84.980726049613381,-86.046245928695086,-87.293554926917366,-88.73059554767967,-90.366319988448865,-92.210859999212261,-94.275672585762052,-96.573681754817997,-99.346220674128304,-102.59107502835546,-106.34835522505755,-110.66524491853907,-115.59702920590982,-121.2082845984838,-127.57427716345563,-134.7826158253942,-143.08613837381597,-152.61349859978577,-163.53533804699072,-176.05364680103159,-190.40773685720475,-206.88151232365976,-225.81233497881965,-247.60185754851531,-270.8280134757353,-295.39306036698792,-320.9736017783805,-347.48291788793136,-374.78872559519408,-402.32230016018667,-429.86051057671136,-457.12661911266974,-483.3151974431255,-508.08775588089208,-529.90718372759216,-547.76486917322211,-559.91052912474822,-565.20077790241112,-565.16664015187303,-562.487823222423,-551.36130295096098,-517.88983099884297,-471.44984669129548,-414.98168128679742,-350.73241978005467,-282.37343654194848,-213.24518142066248,-148.23892991292814,-92.653716377774899,-48.385695215691669,-14.628362513647645)
$row5vals = @(-80,-80.221938104380499,-80.441292253816911,-80.727955183422679,-81.089242484154312,-81.528667220547618,-82.048742471251657,-82.651701066716882,-83.33976419140204,-84.115266644570966,-84.980726049613381,-86.046245928695086,-87.293554926917366,-88.73059554767967,-90.366319988448865,-92.210859999212261,-94.275672585762052,-96.573681754817997,-99.346220674128304,-102.59107502835546,-106.34835522505755,-110.66524491853907,-115.59702920590982,-121.2082845984838,-127.57427716345563,-134.7826158253942,-143.08613837381597,-152.61349859978577,-163.53533804699072,-176.05364680103159,-190.40773685720475,-206.88151232365976,-225.81233497881965,-247.60185754851531,-270.8280134757353,-295.39306036698792,-320.9736017783805,-347.48291788793136,-374.78872559519408,-402.32230016018667,-429.86051057671136,-457.12661911266974,-483.3151974431255,-508.08775588089208,-529.90718372759216,-547.76486917322211,-559.91052912474822,-565.20077790241112,-565.16664015187303,-562.487823222423,-551.36130295096098,-517.88983099884297,-471.44984669129548,-414.98168128679742,-350.73241978005467,-282.37343654194848,-213.24518142066248,-148.23892991292814,-92.653716377774899,-48.385695215691669,-14.628362513647645)
$row6vals = @(2210.5500000000002,0,0,1989.4950000000001,0,0,1790.5455000000002,0,0,1611.4909500000001,0,0,1450.3418550000001,0,0,1305.3076695000002,0,0,1174.7769025500002,0,0,1057.2992122950002,0,0,951.56929106550024,0,0,856.41236195895021,0,0,770.77112576305524,0,0,693.69401318674977,0,0,624.32461186807484,0,0,555.64890456258661,0,0,483.41454696945033,0,0,406.06821945433825,0,0,324.85457556347063,0,0,243.64093167260296,0,0,168.11224285409602,0,0,104.22959056953954,0,0,56.283978907551358)
$row7vals = @(2210.5500000000002,0,0,2063.2215160144488,0,0,1882.8080036188367,0,0,1709.2914496442083,0,0,1548.107757870913,0,0,1400.2266747581864,0,0,1265.3600787277182,0,0,1142.6680434267498,0,0,1031.1590303291489,0,0,930.12425790876898,0,0,838.69960468102829,0,0,755.89487156573273,0,0,681.03269727904706,0,0,606.47057787820472,0,0,527.46046281648682,0,0,441.49550282763306,0,0,342.98899321869982,0,0,269.66796579728577,0,0,189.73508154216125,0,0,119.40335697364351,0,0,65.653826382205708)
$row8vals = @(2210.5500000000002,0,0,1936.0518028052484,0,0,1721.4454032739382,0,0,1536.7528388427074,0,0,1374.5375449647067,0,0,1230.8726159976718,0,0,1103.0807329563593,0,0,989.21394837205116,0,0,887.70944819490865,0,0,796.97075391664043,0,0,715.76055114107919,0,0,643.16297729171379,0,0,578.14594308184928,0,0,514.21020439099971,0,0,447.52622922857421,0,0,377.43031433438631,0,0,311.14322108830936,0,0,221.75414952256722,0,0,149.14518180503313,0,0,90.491507554494788,0,0,47.475045860823556)
$row9vals = @(2210.5500000000002,0,0,1936.0518028052484,0,0,1721.4454032739382,0,0,1536.7528388427074,0,0,1374.5375449647067,0,0,1230.8726159976718,0,0,1103.0807329563593,0,0,989.21394837205116,0,0,887.70944819490865,0,0,796.97075391664043,0,0,715.76055114107919,0,0,643.16297729171379,0,0,578.14594308184928,0,0,514.21020439099971,0,0,447.52622922857421,0,0,377.43031433438631,0,0,311.14322108830936,0,0,221.75414952256722,0,0,149.14518180503313,0,0,90.491507554494788,0,0,47.475045860823556)

for ($i = 0; $i -lt $row2vals.Length; $i++) { $ws.Cells.Item(2, 6+$i).Value = $row2vals[$i] }
for ($i = 0; $i -lt $row3vals.Length; $i++) { $ws.Cells.Item(3, 6+$i).Value = $row3vals[$i] }
for ($i = 0; $i -lt $row4vals.Length; $i++) { $ws.Cells.Item(4, 6+$i).Value = $row4vals[$i] }
for ($i = 0; $i -lt $row5vals.Length; $i++) { $ws.Cells.Item(5, 6+$i).Value = $row5vals[$i] }
for ($i = 0; $i -lt $row6vals.Length; $i++) { $ws.Cells.Item(6, 6+$i).Value = $row6vals[$i] }
for ($i = 0; $i -lt $row7vals.Length; $i++) { $ws.Cells.Item(7, 6+$i).Value = $row7vals[$i] }
for ($i = 0; $i -lt $row8vals.Length; $i++) { $ws.Cells.Item(8, 6+$i).Value = $row8vals[$i] }
for ($i = 0; $i -lt $row9vals.Length; $i++) { $ws.Cells.Item(9, 6+$i).Value = $row9vals[$i] }

# --- Resize the Cashflow table and its data validations to cover the new rows ---
$tbl = $ws.ListObjects.Item("Table_Cashflow")
$tbl.Resize($ws.Range("A1:BN9"))

$ws.Range("E2:E9").Validation.Delete()
$ws.Range("E2:E9").Validation.Add(3, 1, 1, "Novelty_SystemName")
$ws.Range("D2:D9").Validation.Delete()
$ws.Range("D2:D9").Validation.Add(3, 1, 1, "VariableType_SystemName")
$ws.Range("B2:C9").Validation.Delete()
$ws.Range("B2:C9").Validation.Add(3, 1, 1, "AmountType_SystemName")
$ws.Range("A2:A9").Validation.Delete()
$ws.Range("A2:A9").Validation.Add(3, 1, 1, "GroupOfContract_SystemName")

# --- Drop the custom row heights on rows 1 & 2 (back to standard height) ---
$ws.Range("1:2").EntireRow.AutoFit()

# --- Make Cashflow the active sheet / tab and set its scroll selection ---
$ws.Activate()
$ws.Range("H17").Select()
